# The commit refactors the pathway/SWR analysis pipeline: a new results
# column is inserted before the existing column B, shifting the old
# B/C/D columns to C/D/E. The new column B is populated with placeholder
# zeros (the header cell B1 keeps the bold/bordered/centered style shared
# by the other header cells, while B2:B4 stay plain like the other data
# cells), and the last column (E) is refreshed with updated pipeline
# values, with the former last row's value dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this shifts old B->C, C->D, D->E and updates
# the dimension/row spans automatically.
$ws.Columns("B").Insert()

# The insert carries column A's formatting onto the newly created
# B2:B4 cells - clear that back to the default (unstyled) look used by
# the other plain data cells.
$ws.Range("B2:B4").ClearFormats()

# New header cell B1 takes on the same style as the other header cells
# (bold font, thin border, centered/top aligned) by copying the format
# from a neighbouring header cell, then set its value.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = 0

# New data column B is all zeros.
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0

# Refreshed pipeline results land in the (shifted) last column E; row 4's
# value is dropped entirely as part of the refactor.
$ws.Range("E2").Value = -106017.1405710956
$ws.Range("E3").Value = -106016.8796549928
$ws.Range("E4").ClearContents()
